# rekap_tugas.xlsx update
# - Adek (row 3) finished the assignment: Status -> "Sudah", Link File -> submission URL
# - Two new students appended: Budi (row 4), Raka Tegar W (row 5), both "Belum" / "-"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adek's row: mark as submitted, with the link to the uploaded file
$ws.Range("C3").Value = "Sudah"
$ws.Range("D3").Value = '"https://wgdxgzraacfhfbxvxuzy.supabase.co";/storage/v1/object/public/submissions/submissions/1742193980866.pdf'

# New row for Budi
$ws.Range("A4").Value = "Budi"
# Phone number must stay text (leading zeros/long digit strings), not be coerced to a number
$ws.Range("B4").Value = "'62895378394026"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "Belum"
$ws.Range("D4").Value = "-"

# New row for Raka Tegar W
$ws.Range("A5").Value = "Raka Tegar W"
$ws.Range("B5").Value = "'62895396334563"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "Belum"
$ws.Range("D5").Value = "-"
